# Update with restock suggestion
# Applies the "Forecast Comparison" sheet restock-suggestion refresh:
#  - fills in Week_Start_Date (col B) for every forecast week
#  - recomputes Inventory Coverage (L), Stockout Risk (M), Reorder Urgency (N)
#    and Seasonality Index (P)
#  - drops the "Sales Volume Rank" column (Q), shifting "Lifecycle Stage" left
#  - refreshes the Summary sheet's Max/Min Forecast Week cells to "N/A"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Week_Start_Date (col B) values look like dates, but must stay literal text,
# not get auto-converted to a date serial number. Force text entry by
# switching the cell to a text format for the write, then restore General.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2025-02-02"
$ws.Range("B2").NumberFormat = "General"
$ws.Range("L2").Value = 8.33
$ws.Range("M2").Value = "Low"
$ws.Range("N2").Value = "Normal"
$ws.Range("P2").Value = 0.87

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2025-02-09"
$ws.Range("B3").NumberFormat = "General"
$ws.Range("L3").Value = 6.94
$ws.Range("M3").Value = "Low"
$ws.Range("N3").Value = "Normal"
$ws.Range("P3").Value = 0.92

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "2025-02-16"
$ws.Range("B4").NumberFormat = "General"
$ws.Range("L4").Value = 5.68
$ws.Range("M4").Value = "Low"
$ws.Range("N4").Value = "Normal"
$ws.Range("P4").Value = 0.85

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "2025-02-23"
$ws.Range("B5").NumberFormat = "General"
$ws.Range("L5").Value = 4.65
$ws.Range("M5").Value = "Low"
$ws.Range("N5").Value = "Normal"
$ws.Range("P5").Value = 0.8100000000000001

$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "2025-03-02"
$ws.Range("B6").NumberFormat = "General"
$ws.Range("L6").Value = 3.76
$ws.Range("M6").Value = "Low"
$ws.Range("N6").Value = "Normal"
$ws.Range("P6").Value = 1.12

$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "2025-03-09"
$ws.Range("B7").NumberFormat = "General"
$ws.Range("L7").Value = 2.88
$ws.Range("M7").Value = "Low"
$ws.Range("N7").Value = "Normal"
$ws.Range("P7").Value = 0.95

$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "2025-03-16"
$ws.Range("B8").NumberFormat = "General"
$ws.Range("L8").Value = 1.92
$ws.Range("M8").Value = "Low"
$ws.Range("N8").Value = "Normal"
$ws.Range("P8").Value = 0.9

$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "2025-03-23"
$ws.Range("B9").NumberFormat = "General"
$ws.Range("L9").Value = 0.91
$ws.Range("M9").Value = "Low"
$ws.Range("N9").Value = "Urgent"
$ws.Range("P9").Value = 1.07

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "2025-03-30"
$ws.Range("B10").NumberFormat = "General"
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = "High"
$ws.Range("N10").Value = "Urgent"
$ws.Range("P10").Value = 1.11

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "2025-04-06"
$ws.Range("B11").NumberFormat = "General"
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = "High"
$ws.Range("N11").Value = "Urgent"
$ws.Range("P11").Value = 0.85

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "2025-04-13"
$ws.Range("B12").NumberFormat = "General"
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = "High"
$ws.Range("N12").Value = "Urgent"
$ws.Range("P12").Value = 0.98

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "2025-04-20"
$ws.Range("B13").NumberFormat = "General"
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = "High"
$ws.Range("N13").Value = "Urgent"
$ws.Range("P13").Value = 0.84

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "2025-04-27"
$ws.Range("B14").NumberFormat = "General"
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = "High"
$ws.Range("N14").Value = "Urgent"
$ws.Range("P14").Value = 1.01

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "2025-05-04"
$ws.Range("B15").NumberFormat = "General"
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = "High"
$ws.Range("N15").Value = "Urgent"
$ws.Range("P15").Value = 0.9399999999999999

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "2025-05-11"
$ws.Range("B16").NumberFormat = "General"
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = "High"
$ws.Range("N16").Value = "Urgent"
$ws.Range("P16").Value = 0.89

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "2025-05-18"
$ws.Range("B17").NumberFormat = "General"
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = "High"
$ws.Range("N17").Value = "Urgent"
$ws.Range("P17").Value = 1.03

# "Sales Volume Rank" is no longer tracked; deleting column Q shifts
# "Lifecycle Stage" (formerly R) one column to the left, into Q.
$ws.Columns("Q").Delete()

# Summary sheet: with no sufficient week-over-week variance, the min/max
# forecast week callouts are no longer meaningful.
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B13").Value = "N/A"
$summary.Range("B15").Value = "N/A"
